$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $CellRef, $Val)
    $rng = $Worksheet.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "26.499.12"
Set-TextValue $ws "E2" "  +0.14%  "
Set-TextValue $ws "D3" "1.836.97"
Set-TextValue $ws "E3" "  -0.31%  "
Set-TextValue $ws "E4" "  +0.13%  "
Set-TextValue $ws "D5" "257.32"
Set-TextValue $ws "E5" "  -1.55%  "
Set-TextValue $ws "E6" "  +0.10%  "
Set-TextValue $ws "D7" "0.5240"
Set-TextValue $ws "E7" "  +0.66%  "
Set-TextValue $ws "D8" "0.3138"
Set-TextValue $ws "E8" "  -4.24%  "
Set-TextValue $ws "E9" "  -0.01%  "
Set-TextValue $ws "D10" "18.58"
Set-TextValue $ws "E10" "  -0.51%  "
Set-TextValue $ws "D11" "0.7731"
Set-TextValue $ws "E11" "  +0.04%  "
Set-TextValue $ws "D12" "0.07755"
Set-TextValue $ws "E12" "  +0.73%  "
Set-TextValue $ws "D13" "1.832.34"
Set-TextValue $ws "E13" "  +0.84%  "
Set-TextValue $ws "D14" "87.54"
Set-TextValue $ws "E14" "  -0.79%  "
Set-TextValue $ws "D15" "4.993"
Set-TextValue $ws "E15" "  -0.72%  "
Set-TextValue $ws "E16" "  +0.16%  "
Set-TextValue $ws "D17" "13.77"
Set-TextValue $ws "E17" "  -1.22%  "
Set-TextValue $ws "E18" "  +0.12%  "
Set-TextValue $ws "D19" "0.000007911"
Set-TextValue $ws "E19" "  -0.59%  "
Set-TextValue $ws "D20" "26.516.71"
Set-TextValue $ws "E20" "  +0.36%  "
Set-TextValue $ws "D21" "2.064.69"
Set-TextValue $ws "E21" "  +0.35%  "
Set-TextValue $ws "D22" "4.578"
Set-TextValue $ws "E22" "  +0.00%  "
Set-TextValue $ws "D23" "5.943"
Set-TextValue $ws "E23" "  -0.46%  "
Set-TextValue $ws "D24" "9.296"
Set-TextValue $ws "E24" "  -2.04%  "
Set-TextValue $ws "D25" "142.21"
Set-TextValue $ws "E25" "  -1.40%  "
Set-TextValue $ws "D26" "2.205"
Set-TextValue $ws "E26" "  -0.16%  "
Set-TextValue $ws "D27" "1.670"
Set-TextValue $ws "E27" "  +1.31%  "
Set-TextValue $ws "E28" "  -0.77%  "
Set-TextValue $ws "D29" "110.54"
Set-TextValue $ws "E29" "  -0.76%  "
Set-TextValue $ws "E30" "  -1.33%  "
Set-TextValue $ws "D31" "0.08711"
Set-TextValue $ws "E31" "  -0.18%  "
Set-TextValue $ws "E32" "  -2.10%  "
Set-TextValue $ws "D33" "0.04841"
Set-TextValue $ws "E33" "  +0.78%  "
Set-TextValue $ws "D34" "1.131"
Set-TextValue $ws "E34" "  -0.01%  "
Set-TextValue $ws "D35" "0.7148"
Set-TextValue $ws "E35" "  +0.75%  "
Set-TextValue $ws "E36" "  +0.83%  "
Set-TextValue $ws "D37" "3.078"
Set-TextValue $ws "E37" "  -0.23%  "
Set-TextValue $ws "D38" "2.212"
Set-TextValue $ws "E38" "  -0.99%  "
Set-TextValue $ws "D39" "0.01720"
Set-TextValue $ws "E39" "  -2.38%  "
Set-TextValue $ws "D40" "0.4780"
Set-TextValue $ws "E40" "  -1.42%  "
Set-TextValue $ws "D41" "0.8909"
Set-TextValue $ws "E41" "  -0.47%  "
Set-TextValue $ws "D42" "109.47"
Set-TextValue $ws "E42" "  -1.79%  "
Set-TextValue $ws "D43" "5.906"
Set-TextValue $ws "E43" "  -2.64%  "
Set-TextValue $ws "E44" "  +0.17%  "
Set-TextValue $ws "D45" "7.606"
Set-TextValue $ws "E45" "  -1.70%  "
Set-TextValue $ws "E46" "  -0.92%  "
Set-TextValue $ws "D47" "8.957"
Set-TextValue $ws "E47" "  -0.51%  "
Set-TextValue $ws "B48" "Cronos"
Set-TextValue $ws "C48" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D48" "0.05802"
Set-TextValue $ws "E48" "  -1.17%  "
Set-TextValue $ws "B49" "Algorand"
Set-TextValue $ws "C49" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws "D49" "0.1224"
Set-TextValue $ws "E49" "  +0.42%  "
Set-TextValue $ws "D50" "34.60"
Set-TextValue $ws "E50" "  -1.20%  "
Set-TextValue $ws "D51" "0.8903"
Set-TextValue $ws "E51" "  +0.24%  "
